$d = $word.ActiveDocument

$replacements = @(
    @{old="11×86=946"; new="95×23=2185"},
    @{old="50×14=700"; new="43×68=2924"},
    @{old="63×23=1449"; new="29×22=638"},
    @{old="59×17=1003"; new="60×66=3960"},
    @{old="71×18=1278"; new="48×33=1584"},
    @{old="73×81=5913"; new="19×43=817"},
    @{old="19×46=874"; new="81×99=8019"},
    @{old="68×30=2040"; new="85×52=4420"},
    @{old="94×31=2914"; new="72×68=4896"},
    @{old="39×55=2145"; new="70×75=5250"},
    @{old="78×50=3900"; new="80×97=7760"},
    @{old="76×55=4180"; new="37×13=481"},
    @{old="24×43=1032"; new="13×70=910"},
    @{old="22×14=308"; new="32×77=2464"},
    @{old="59×32=1888"; new="17×92=1564"},
    @{old="56×25=1400"; new="29×32=928"},
    @{old="60×12=720"; new="43×41=1763"},
    @{old="88×17=1496"; new="40×41=1640"},
    @{old="20×63=1260"; new="75×54=4050"},
    @{old="11×93=1023"; new="38×50=1900"},
    @{old="40×75=3000"; new="60×57=3420"},
    @{old="24×54=1296"; new="89×60=5340"},
    @{old="42×84=3528"; new="95×15=1425"},
    @{old="89×12=1068"; new="26×34=884"},
    @{old="64×52=3328"; new="67×64=4288"}
)

foreach ($pair in $replacements) {
    $d.Content.Find.Execute($pair.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $pair.new, 2)
}
